$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the A:E contents of row pairs (2016年B<->2016年C, 2017年B<->2017年C,
# 2018年B<->2018年C, 2019年B<->2019年C) so that the "C" period sorts before
# the "B" period within each year block.
$rowPairs = @(
    @(3, 4),
    @(7, 8),
    @(11, 12),
    @(15, 16)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 1; $col -le 5; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $v1 = $c1.Value()
        $v2 = $c2.Value()
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# Remove the now-unwanted "合成洗涤剂产销率" (F) and "合成洗涤剂销售量" (G)
# columns entirely, shrinking the used range back down to A1:E17.
$ws.Range("F1:G17").Delete()
